# Updates the cryptocurrency price/volume table on Sheet1 to reflect the
# latest scrape (GitHub Actions data refresh).
#
# Column D ("Price") and most of column E ("Volume(1h)") are stored as
# literal text in the workbook (prices use "."-grouped formatting that
# isn't valid numeric literal in several rows, e.g. "29.205.54"), so the
# text is written back verbatim. For the handful of Price cells whose new
# text *would* parse as a plain number (e.g. "1.003"), a leading apostrophe
# forces Excel to keep the cell as text instead of silently converting it
# to a numeric value - matching how the source data is represented.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '29.205.54'
$ws.Range("E2").Value = '  -0.66%  '
$ws.Range("D3").Value = '1.825.63'
$ws.Range("E3").Value = '  -0.87%  '
$ws.Range("D4").Value = '''1.000'
$ws.Range("E4").Value = '  +0.09%  '
$ws.Range("D5").Value = '''234.23'
$ws.Range("E5").Value = '  -2.12%  '
$ws.Range("D6").Value = '''0.6003'
$ws.Range("E6").Value = '  -4.15%  '
$ws.Range("D7").Value = '''1.003'
$ws.Range("E7").Value = '  +0.21%  '
$ws.Range("D8").Value = '''0.07019'
$ws.Range("E8").Value = '  -5.66%  '
$ws.Range("D9").Value = '''0.2778'
$ws.Range("E9").Value = '  -4.04%  '
$ws.Range("D10").Value = '''23.38'
$ws.Range("E10").Value = '  -6.50%  '
$ws.Range("D12").Value = '1.829.10'
$ws.Range("E12").Value = '  -0.72%  '
$ws.Range("E13").Value = '  -4.00%  '
$ws.Range("D14").Value = '''0.000009925'
$ws.Range("E14").Value = '  -3.62%  '
$ws.Range("D15").Value = '''0.6237'
$ws.Range("E15").Value = '  -7.76%  '
$ws.Range("D16").Value = '''78.74'
$ws.Range("E16").Value = '  -3.79%  '
$ws.Range("D17").Value = '29.209.07'
$ws.Range("E17").Value = '  -0.76%  '
$ws.Range("D18").Value = '''5.815'
$ws.Range("E18").Value = '  -6.70%  '
$ws.Range("D19").Value = '''223.25'
$ws.Range("E19").Value = '  -3.97%  '
$ws.Range("E20").Value = '  +0.20%  '
$ws.Range("D21").Value = '''11.64'
$ws.Range("E21").Value = '  -5.57%  '
$ws.Range("D22").Value = '''6.980'
$ws.Range("E22").Value = '  -4.74%  '
$ws.Range("D23").Value = '''1.003'
$ws.Range("E23").Value = '  +0.24%  '
$ws.Range("D24").Value = '''155.31'
$ws.Range("E24").Value = '  -1.71%  '
$ws.Range("D25").Value = '''7.957'
$ws.Range("E25").Value = '  -6.25%  '
$ws.Range("D26").Value = '''0.1291'
$ws.Range("E26").Value = '  -4.41%  '
$ws.Range("E27").Value = '  -5.05%  '
$ws.Range("D28").Value = '''1.482'
$ws.Range("E28").Value = '  +0.80%  '
$ws.Range("D29").Value = '''0.06225'
$ws.Range("E29").Value = '  -12.95%  '
$ws.Range("D30").Value = '''1.438'
$ws.Range("E30").Value = '  -2.99%  '
$ws.Range("D31").Value = '''3.821'
$ws.Range("E31").Value = '  -5.27%  '
$ws.Range("D32").Value = '''3.777'
$ws.Range("E32").Value = '  -6.58%  '
$ws.Range("D33").Value = '''1.107'
$ws.Range("E33").Value = '  -2.91%  '
$ws.Range("D34").Value = '''1.734'
$ws.Range("E34").Value = '  -4.82%  '
$ws.Range("D35").Value = '''0.6423'
$ws.Range("E35").Value = '  -8.07%  '
$ws.Range("D36").Value = '''2.542'
$ws.Range("E36").Value = '  -1.37%  '
$ws.Range("D37").Value = '1.219.66'
$ws.Range("E37").Value = '  -1.20%  '
$ws.Range("E38").Value = '  -3.02%  '
$ws.Range("D39").Value = '''6.522'
$ws.Range("E39").Value = '  -5.72%  '
$ws.Range("D40").Value = '''0.01724'
$ws.Range("E40").Value = '  -6.40%  '
$ws.Range("D41").Value = '''0.8960'
$ws.Range("E41").Value = '  -6.93%  '
$ws.Range("D42").Value = '''1.004'
$ws.Range("E42").Value = '  +0.33%  '
$ws.Range("D43").Value = '1.980.32'
$ws.Range("E43").Value = '  -1.74%  '
$ws.Range("D44").Value = '''99.87'
$ws.Range("E44").Value = '  -1.10%  '
$ws.Range("D45").Value = '''62.25'
$ws.Range("E45").Value = '  -5.01%  '
$ws.Range("D46").Value = '''0.00000000117'
$ws.Range("E46").Value = '  -2.07%  '
$ws.Range("D47").Value = '''8.518'
$ws.Range("E47").Value = '  -4.64%  '
$ws.Range("D48").Value = '''0.4550'
$ws.Range("E48").Value = '  -0.69%  '
$ws.Range("B49").Value = 'RenderToken'
$ws.Range("C49").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D49").Value = '''1.568'
$ws.Range("E49").Value = '  -9.16%  '
$ws.Range("B50").Value = 'Cronos'
$ws.Range("C50").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D50").Value = '''0.05484'
$ws.Range("E50").Value = '  -3.05%  '
$ws.Range("D51").Value = '''6.397'
$ws.Range("E51").Value = '  -8.18%  '

# Rows 49/50 swapped which coin occupies each row (Cronos <-> RenderToken)
# along with their Link/Price/Volume figures; handled above via direct
# cell writes rather than a row move.
